$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date value for all data rows (2-416)
# from 45172 (2023-09-03) to 45175 (2023-09-06)
$ws.Range("C2:C416").Value = 45175
